$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy formatting from column G to the new column H for rows 1-8, while G
#    still carries its original look (so H ends up identical to G, as in the
#    source file before the table grew by one column).
# ---------------------------------------------------------------------------
$ws.Range("G1:G8").Copy()
$ws.Range("H1:H8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Grow the merged header cell from F1:G1 to F1:H1.  The engine's Merge
#    recomputes an "outside border" for the whole span, which clobbers the
#    bespoke per-cell look the template relies on (F1 = full box, G1 = top+
#    bottom only, H1/former-G1 = top+bottom+right).  So merge first, then
#    restore the exact original look on every involved cell afterwards.
# ---------------------------------------------------------------------------
$ws.Range("F1:G1").UnMerge()
$ws.Range("F1:H1").MergeCells = $true

# Restore F1's original "full box" look by pulling it back from D1, which is
# untouched and still uses the very same style.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Restore G1 (now an interior cell of the merge: top+bottom border only) and
# H1 (now the rightmost cell: top+bottom+right border) explicitly.
$green = 5880731

$ws.Range("G1").ClearFormats()
$ws.Range("G1").Borders.Item(8).LineStyle = 1
$ws.Range("G1").Borders.Item(8).Color = $green
$ws.Range("G1").Borders.Item(9).LineStyle = 1
$ws.Range("G1").Borders.Item(9).Color = $green

$ws.Range("H1").ClearFormats()
$ws.Range("H1").Borders.Item(8).LineStyle = 1
$ws.Range("H1").Borders.Item(8).Color = $green
$ws.Range("H1").Borders.Item(9).LineStyle = 1
$ws.Range("H1").Borders.Item(9).Color = $green
$ws.Range("H1").Borders.Item(10).LineStyle = 1
$ws.Range("H1").Borders.Item(10).Color = $green

# ---------------------------------------------------------------------------
# 3. Updated dimensions (2500->2530 / 580->610) for rows 5-8.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 2530
$ws.Range("C5").Value = 610
$ws.Range("B6").Value = 2530
$ws.Range("C6").Value = 610
$ws.Range("B7").Value = 2530
$ws.Range("C7").Value = 610
$ws.Range("B8").Value = 2530
$ws.Range("C8").Value = 610

# ---------------------------------------------------------------------------
# 4. Column width adjustments.  The engine's ColumnWidth setter re-adds the
#    standard ~0.8333 character padding on export, so the value written here
#    is shifted down by 5/6 to land exactly on the intended raw XML widths
#    (8 and 13).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666

# ---------------------------------------------------------------------------
# 5. Extend the table ("MyTable") to include the new column, then set the new
#    header label - the ListColumn name syncs from the header cell text once
#    the table range covers it, so the header text must be written *after*
#    Resize() for the table XML to pick up the real name instead of the
#    default placeholder.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:H8"))
$ws.Range("H2").Value = "Tekst oper.10"

Write-Host "Edit applied: added column H (Tekst oper.10) to table and sheet"
